$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Detach the "_GoBack" bookmark from its current position (right
#    after "Common", at the end of paragraph 8) so it can be re-added
#    at the end of the document once the new paragraphs exist.
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# 2. Append two new paragraphs after "Паджинация" (currently the last
#    paragraph in the body, right before sectPr):
#       - "Hidden/" + "ложное удаление"
#       - "Асинхронность"
# ---------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)

# --- paragraph: Hidden/ложное удаление ---------------------------------
$endRng = $lastPara.Range.Duplicate()
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

$newParaIndex = $lastParaIndex + 1
$hiddenPara = $d.Paragraphs.Item($newParaIndex)
$insPoint = $hiddenPara.Range.Duplicate()
$insPoint.Collapse(0)
$insPoint.InsertAfter("ложное удаление")

# Paste an unformatted ("no rPr") run in front of the text just typed,
# then overwrite its placeholder text with "Hidden/" so the run keeps
# the plain (non ru-RU) formatting that the source run ("CRUD") has.
# NOTE: Paste() must target a Range built directly via $d.Range(x, x)
# (collapsed at construction) -- a Range.Duplicate() + Collapse() still
# carries its pre-collapse extent into Paste() and clobbers content.
$crudRun = $d.Range($d.Paragraphs.Item(10).Range.Start, $d.Paragraphs.Item(10).Range.End - 1)
$crudRun.Copy()

$hiddenParaStart = $hiddenPara.Range.Start
$hiddenStart = $d.Range($hiddenParaStart, $hiddenParaStart)
$hiddenStart.Paste()

$plainRun = $d.Range($hiddenParaStart, $hiddenParaStart + 4)
$plainRun.Text = "Hidden/"

# --- paragraph: Асинхронность ------------------------------------------
$hiddenParaNow = $d.Paragraphs.Item($newParaIndex)
$endRng2 = $hiddenParaNow.Range.Duplicate()
$endRng2.Collapse(0)
$endRng2.InsertParagraphAfter()

$asyncParaIndex = $newParaIndex + 1
$asyncPara = $d.Paragraphs.Item($asyncParaIndex)
$insPoint2 = $asyncPara.Range.Duplicate()
$insPoint2.Collapse(0)
$insPoint2.InsertAfter("Асинхронность")

# ---------------------------------------------------------------------
# 3. Re-create the "_GoBack" bookmark, collapsed, right after the new
#    "Асинхронность" text (before the paragraph mark) -- mirroring
#    where it originally sat after "Common".
# ---------------------------------------------------------------------
$asyncParaNow = $d.Paragraphs.Item($asyncParaIndex)
$bmSpot = $asyncParaNow.Range.Duplicate()
$bmSpot.Collapse(0)
$bmSpot.MoveEnd(1, -1)
$bmSpot.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $bmSpot)
$placeholder = $d.Range($bmSpot.Start, $bmSpot.End)
$placeholder.Text = ""
